# Update cryptocurrency price/volume data per upstream refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.528.63"
$ws.Range("E2").Value = "  -0.10%  "
$ws.Range("D3").Value = "1.618.50"
$ws.Range("E3").Value = "  -1.32%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").Value = "'210.82"
$ws.Range("E5").Value = "  -0.79%  "
$ws.Range("E6").Value = "  -1.95%  "
$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("E8").Value = "  -0.72%  "
$ws.Range("E9").Value = "  +2.35%  "
$ws.Range("E10").Value = "  +0.15%  "
$ws.Range("D11").Value = "'0.0885"
$ws.Range("E11").Value = "  -0.37%  "
$ws.Range("D12").Value = "1.845.98"
$ws.Range("E12").Value = "  -1.47%  "
$ws.Range("D13").Value = "1.620.60"
$ws.Range("E13").Value = "  -1.42%  "
$ws.Range("E14").Value = "  -0.08%  "
$ws.Range("E15").Value = "  -2.18%  "
$ws.Range("D16").Value = "'65.03"
$ws.Range("E16").Value = "  +1.71%  "
$ws.Range("D17").Value = "27.511.14"
$ws.Range("E17").Value = "  -0.07%  "
$ws.Range("D18").Value = "'229.95"
$ws.Range("E18").Value = "  +1.08%  "
$ws.Range("E19").Value = "  -0.57%  "
$ws.Range("D20").Value = "'7.53"
$ws.Range("E20").Value = "  -0.85%  "
$ws.Range("E21").Value = "  +0.13%  "
$ws.Range("D22").Value = "'4.29"
$ws.Range("E22").Value = "  +0.04%  "
$ws.Range("D23").Value = "'10.12"
$ws.Range("E23").Value = "  +1.19%  "
$ws.Range("E24").Value = "  +7.35%  "
$ws.Range("D25").Value = "'149.38"
$ws.Range("E25").Value = "  +0.14%  "
$ws.Range("E26").Value = "  -1.10%  "
$ws.Range("B27").Value = "Cosmos"
$ws.Range("C27").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D27").Value = "'6.81"
$ws.Range("E27").Value = "  -1.98%  "
$ws.Range("B28").Value = "BinanceUSD"
$ws.Range("C28").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D28").Value = "'1.00"
$ws.Range("E28").Value = "  +0.05%  "
$ws.Range("D29").Value = "'15.55"
$ws.Range("E29").Value = "  -0.16%  "
$ws.Range("E30").Value = "  -0.55%  "
$ws.Range("D31").Value = "'0.0482"
$ws.Range("E31").Value = "  -0.70%  "
$ws.Range("E32").Value = "  -0.88%  "
$ws.Range("D33").Value = "1.442.38"
$ws.Range("E33").Value = "  +1.12%  "
$ws.Range("E34").Value = "  -3.46%  "
$ws.Range("E35").Value = "  -3.27%  "
$ws.Range("E36").Value = "  -0.23%  "
$ws.Range("D37").Value = "'0.935"
$ws.Range("E37").Value = "  +3.87%  "
$ws.Range("D38").Value = "'0.560"
$ws.Range("E38").Value = "  -2.22%  "
$ws.Range("E39").Value = "  +0.19%  "
$ws.Range("D40").Value = "'0.860"
$ws.Range("E40").Value = "  -1.62%  "
$ws.Range("D41").Value = "'69.15"
$ws.Range("E41").Value = "  +6.26%  "
$ws.Range("E42").Value = "  +0.02%  "
$ws.Range("D43").Value = "'0.996"
$ws.Range("E43").Value = "  -3.09%  "
$ws.Range("E44").Value = "  +0.11%  "
$ws.Range("E45").Value = "  -1.93%  "
$ws.Range("E46").Value = "  -1.94%  "
$ws.Range("D47").Value = "1.757.14"
$ws.Range("E47").Value = "  -1.49%  "
$ws.Range("E48").Value = "  +0.05%  "
$ws.Range("D49").Value = "'86.32"
$ws.Range("E49").Value = "  +0.08%  "
$ws.Range("E50").Value = "  -1.37%  "
$ws.Range("D51").Value = "'0.0994"
$ws.Range("E51").Value = "  +1.12%  "
